$d = $word.ActiveDocument

# Update the date/day line in the first paragraph.
$d.Content.Find.Execute("2024-04-06 Saturday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2024-04-07 Sunday", 2)

# Update the division problems in the table, addressed by (row, column) so
# that the two cells that happen to share the original text ("64÷8=") can
# each be mapped to their own distinct replacement.
$table = $d.Tables.Item(1)

$updates = @(
    @(1,1,'36÷8='),
    @(1,2,'15÷2='),
    @(1,3,'51÷2='),
    @(1,4,'39÷5='),
    @(1,5,'76÷8='),

    @(5,1,'55÷7='),
    @(5,2,'26÷9='),
    @(5,3,'98÷2='),
    @(5,4,'72÷5='),
    @(5,5,'64÷4='),

    @(9,1,'98÷7='),
    @(9,2,'63÷5='),
    @(9,3,'60÷2='),
    @(9,4,'27÷5='),
    @(9,5,'22÷9='),

    @(13,1,'45÷9='),
    @(13,2,'10÷7='),
    @(13,3,'98÷3='),
    @(13,4,'23÷9='),
    @(13,5,'90÷3='),

    @(17,1,'81÷2='),
    @(17,2,'13÷8='),
    @(17,3,'90÷4='),
    @(17,4,'23÷9='),
    @(17,5,'93÷9=')
)

foreach ($u in $updates) {
    $row = $u[0]
    $col = $u[1]
    $newText = $u[2]
    $cell = $table.Cell($row, $col)
    $cellRange = $cell.Range
    $cellRange.End = $cellRange.End - 1
    $cellRange.Text = $newText
}
